$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 1).Value = '(''looks'', ''like'')'
$ws1.Cells.Item(2, 2).Value = 14
$ws1.Cells.Item(2, 3).Value = 0.03526448362720403
$ws1.Cells.Item(3, 1).Value = '(''year'', ''old'')'
$ws1.Cells.Item(3, 2).Value = 11
$ws1.Cells.Item(3, 3).Value = 0.02770780856423174
$ws1.Cells.Item(4, 1).Value = '(''waste'', ''time'')'
$ws1.Cells.Item(4, 2).Value = 10
$ws1.Cells.Item(4, 3).Value = 0.02518891687657431
$ws1.Cells.Item(5, 1).Value = '(''really'', ''bad'')'
$ws1.Cells.Item(5, 2).Value = 10
$ws1.Cells.Item(5, 3).Value = 0.02518891687657431
$ws1.Cells.Item(6, 1).Value = '(''bad'', ''acting'')'
$ws1.Cells.Item(6, 2).Value = 10
$ws1.Cells.Item(6, 3).Value = 0.02518891687657431
$ws1.Cells.Item(7, 1).Value = '(''one'', ''worst'')'
$ws1.Cells.Item(7, 2).Value = 9
$ws1.Cells.Item(7, 3).Value = 0.02267002518891688
$ws1.Cells.Item(8, 1).Value = '(''even'', ''worse'')'
$ws1.Cells.Item(8, 2).Value = 9
$ws1.Cells.Item(8, 3).Value = 0.02267002518891688
$ws1.Cells.Item(9, 1).Value = '(''stay'', ''away'')'
$ws1.Cells.Item(9, 2).Value = 9
$ws1.Cells.Item(9, 3).Value = 0.02267002518891688
$ws1.Cells.Item(10, 1).Value = '(''movie'', ''ever'')'
$ws1.Cells.Item(10, 2).Value = 9
$ws1.Cells.Item(10, 3).Value = 0.02267002518891688
$ws1.Cells.Item(11, 1).Value = '(''movie'', ''see'')'
$ws1.Cells.Item(11, 2).Value = 8
$ws1.Cells.Item(11, 3).Value = 0.02015113350125945
$ws1.Cells.Item(12, 1).Value = '(''movies'', ''like'')'
$ws1.Cells.Item(12, 2).Value = 8
$ws1.Cells.Item(12, 3).Value = 0.02015113350125945
$ws1.Cells.Item(13, 1).Value = '(''tv'', ''series'')'
$ws1.Cells.Item(13, 2).Value = 8
$ws1.Cells.Item(13, 3).Value = 0.02015113350125945
$ws1.Cells.Item(14, 1).Value = '(''horror'', ''movie'')'
$ws1.Cells.Item(14, 2).Value = 8
$ws1.Cells.Item(14, 3).Value = 0.02015113350125945
$ws1.Cells.Item(15, 1).Value = '(''movie'', ''would'')'
$ws1.Cells.Item(15, 2).Value = 8
$ws1.Cells.Item(15, 3).Value = 0.02015113350125945
$ws1.Cells.Item(16, 1).Value = '(''worth'', ''watching'')'
$ws1.Cells.Item(16, 2).Value = 7
$ws1.Cells.Item(16, 3).Value = 0.01763224181360202
$ws1.Cells.Item(17, 1).Value = '(''high'', ''school'')'
$ws1.Cells.Item(17, 2).Value = 7
$ws1.Cells.Item(17, 3).Value = 0.01763224181360202
$ws1.Cells.Item(18, 1).Value = '(''even'', ''get'')'
$ws1.Cells.Item(18, 2).Value = 7
$ws1.Cells.Item(18, 3).Value = 0.01763224181360202
$ws1.Cells.Item(19, 1).Value = '(''seems'', ''like'')'
$ws1.Cells.Item(19, 2).Value = 7
$ws1.Cells.Item(19, 3).Value = 0.01763224181360202
$ws1.Cells.Item(20, 1).Value = '(''like'', ''movie'')'
$ws1.Cells.Item(20, 2).Value = 7
$ws1.Cells.Item(20, 3).Value = 0.01763224181360202
$ws1.Cells.Item(21, 1).Value = '(''like'', ''one'')'
$ws1.Cells.Item(21, 2).Value = 7
$ws1.Cells.Item(21, 3).Value = 0.01763224181360202
$ws1.Cells.Item(22, 1).Value = '(''movie'', ''also'')'
$ws1.Cells.Item(22, 2).Value = 6
$ws1.Cells.Item(22, 3).Value = 0.01511335012594458
$ws1.Cells.Item(23, 1).Value = '(''movie'', ''could'')'
$ws1.Cells.Item(23, 2).Value = 6
$ws1.Cells.Item(23, 3).Value = 0.01511335012594458
$ws1.Cells.Item(24, 1).Value = '(''worst'', ''movie'')'
$ws1.Cells.Item(24, 2).Value = 6
$ws1.Cells.Item(24, 3).Value = 0.01511335012594458
$ws1.Cells.Item(25, 1).Value = '(''movie'', ''one'')'
$ws1.Cells.Item(25, 2).Value = 6
$ws1.Cells.Item(25, 3).Value = 0.01511335012594458
$ws1.Cells.Item(26, 1).Value = '(''movie'', ''first'')'
$ws1.Cells.Item(26, 2).Value = 6
$ws1.Cells.Item(26, 3).Value = 0.01511335012594458
$ws1.Cells.Item(27, 1).Value = '(''would'', ''better'')'
$ws1.Cells.Item(27, 2).Value = 6
$ws1.Cells.Item(27, 3).Value = 0.01511335012594458
$ws1.Cells.Item(28, 1).Value = '(''first'', ''movie'')'
$ws1.Cells.Item(28, 2).Value = 6
$ws1.Cells.Item(28, 3).Value = 0.01511335012594458
$ws1.Cells.Item(29, 1).Value = '(''sounds'', ''like'')'
$ws1.Cells.Item(29, 2).Value = 6
$ws1.Cells.Item(29, 3).Value = 0.01511335012594458
$ws1.Cells.Item(30, 1).Value = '(''bad'', ''one'')'
$ws1.Cells.Item(30, 2).Value = 6
$ws1.Cells.Item(30, 3).Value = 0.01511335012594458
$ws1.Cells.Item(31, 1).Value = '(''movie'', ''think'')'
$ws1.Cells.Item(31, 2).Value = 6
$ws1.Cells.Item(31, 3).Value = 0.01511335012594458
$ws1.Cells.Item(32, 1).Value = '(''bad'', ''movie'')'
$ws1.Cells.Item(32, 2).Value = 6
$ws1.Cells.Item(32, 3).Value = 0.01511335012594458
$ws1.Cells.Item(33, 1).Value = '(''movie'', ''bad'')'
$ws1.Cells.Item(33, 2).Value = 6
$ws1.Cells.Item(33, 3).Value = 0.01511335012594458
$ws1.Cells.Item(34, 1).Value = '(''martial'', ''arts'')'
$ws1.Cells.Item(34, 2).Value = 5
$ws1.Cells.Item(34, 3).Value = 0.01259445843828715
$ws1.Cells.Item(35, 1).Value = '(''writer'', ''director'')'
$ws1.Cells.Item(35, 2).Value = 5
$ws1.Cells.Item(35, 3).Value = 0.01259445843828715
$ws1.Cells.Item(36, 1).Value = '(''movie'', ''makes'')'
$ws1.Cells.Item(36, 2).Value = 5
$ws1.Cells.Item(36, 3).Value = 0.01259445843828715
$ws1.Cells.Item(37, 1).Value = '(''watching'', ''movie'')'
$ws1.Cells.Item(37, 2).Value = 5
$ws1.Cells.Item(37, 3).Value = 0.01259445843828715
$ws1.Cells.Item(38, 1).Value = '(''film'', ''making'')'
$ws1.Cells.Item(38, 2).Value = 5
$ws1.Cells.Item(38, 3).Value = 0.01259445843828715
$ws1.Cells.Item(39, 1).Value = '(''get'', ''wrong'')'
$ws1.Cells.Item(39, 2).Value = 5
$ws1.Cells.Item(39, 3).Value = 0.01259445843828715
$ws1.Cells.Item(40, 1).Value = '(''bad'', ''bad'')'
$ws1.Cells.Item(40, 2).Value = 5
$ws1.Cells.Item(40, 3).Value = 0.01259445843828715
$ws1.Cells.Item(41, 1).Value = '(''look'', ''like'')'
$ws1.Cells.Item(41, 2).Value = 5
$ws1.Cells.Item(41, 3).Value = 0.01259445843828715
$ws1.Cells.Item(42, 1).Value = '(''movie'', ''even'')'
$ws1.Cells.Item(42, 2).Value = 5
$ws1.Cells.Item(42, 3).Value = 0.01259445843828715
$ws1.Cells.Item(43, 1).Value = '(''movie'', ''like'')'
$ws1.Cells.Item(43, 2).Value = 5
$ws1.Cells.Item(43, 3).Value = 0.01259445843828715
$ws1.Cells.Item(44, 1).Value = '(''slow'', ''moving'')'
$ws1.Cells.Item(44, 2).Value = 5
$ws1.Cells.Item(44, 3).Value = 0.01259445843828715
$ws1.Cells.Item(45, 1).Value = '(''horror'', ''films'')'
$ws1.Cells.Item(45, 2).Value = 5
$ws1.Cells.Item(45, 3).Value = 0.01259445843828715
$ws1.Cells.Item(46, 1).Value = '(''seen'', ''film'')'
$ws1.Cells.Item(46, 2).Value = 5
$ws1.Cells.Item(46, 3).Value = 0.01259445843828715
$ws1.Cells.Item(47, 1).Value = '(''entire'', ''film'')'
$ws1.Cells.Item(47, 2).Value = 5
$ws1.Cells.Item(47, 3).Value = 0.01259445843828715
$ws1.Cells.Item(48, 1).Value = '(''never'', ''seen'')'
$ws1.Cells.Item(48, 2).Value = 5
$ws1.Cells.Item(48, 3).Value = 0.01259445843828715
$ws1.Cells.Item(49, 1).Value = '(''anyone'', ''else'')'
$ws1.Cells.Item(49, 2).Value = 5
$ws1.Cells.Item(49, 3).Value = 0.01259445843828715
$ws1.Cells.Item(50, 1).Value = '(''plot'', ''line'')'
$ws1.Cells.Item(50, 2).Value = 5
$ws1.Cells.Item(50, 3).Value = 0.01259445843828715
$ws1.Cells.Item(51, 1).Value = '(''pretty'', ''bad'')'
$ws1.Cells.Item(51, 2).Value = 5
$ws1.Cells.Item(51, 3).Value = 0.01259445843828715
$ws1.Cells.Item(52, 1).Value = '(''five'', ''minutes'')'
$ws1.Cells.Item(52, 2).Value = 5
$ws1.Cells.Item(52, 3).Value = 0.01259445843828715
$ws1.Cells.Item(53, 1).Value = '(''many'', ''times'')'
$ws1.Cells.Item(53, 2).Value = 5
$ws1.Cells.Item(53, 3).Value = 0.01259445843828715
$ws1.Cells.Item(54, 1).Value = '(''rest'', ''movie'')'
$ws1.Cells.Item(54, 2).Value = 5
$ws1.Cells.Item(54, 3).Value = 0.01259445843828715
$ws1.Cells.Item(55, 1).Value = '(''rest'', ''cast'')'
$ws1.Cells.Item(55, 2).Value = 5
$ws1.Cells.Item(55, 3).Value = 0.01259445843828715
$ws1.Cells.Item(56, 1).Value = '(''none'', ''characters'')'
$ws1.Cells.Item(56, 2).Value = 5
$ws1.Cells.Item(56, 3).Value = 0.01259445843828715
$ws1.Cells.Item(57, 1).Value = '(''mst'', ''k'')'
$ws1.Cells.Item(57, 2).Value = 5
$ws1.Cells.Item(57, 3).Value = 0.01259445843828715
$ws1.Cells.Item(58, 1).Value = '(''bad'', ''movies'')'
$ws1.Cells.Item(58, 2).Value = 4
$ws1.Cells.Item(58, 3).Value = 0.01007556675062972
$ws1.Cells.Item(59, 1).Value = '(''little'', ''boy'')'
$ws1.Cells.Item(59, 2).Value = 4
$ws1.Cells.Item(59, 3).Value = 0.01007556675062972
$ws1.Cells.Item(60, 1).Value = '(''say'', ''something'')'
$ws1.Cells.Item(60, 2).Value = 4
$ws1.Cells.Item(60, 3).Value = 0.01007556675062972
$ws1.Cells.Item(61, 1).Value = '(''film'', ''even'')'
$ws1.Cells.Item(61, 2).Value = 4
$ws1.Cells.Item(61, 3).Value = 0.01007556675062972
$ws1.Cells.Item(62, 1).Value = '(''movie'', ''horrible'')'
$ws1.Cells.Item(62, 2).Value = 4
$ws1.Cells.Item(62, 3).Value = 0.01007556675062972
$ws1.Cells.Item(63, 1).Value = '(''time'', ''money'')'
$ws1.Cells.Item(63, 2).Value = 4
$ws1.Cells.Item(63, 3).Value = 0.01007556675062972
$ws1.Cells.Item(64, 1).Value = '(''george'', ''mildred'')'
$ws1.Cells.Item(64, 2).Value = 4
$ws1.Cells.Item(64, 3).Value = 0.01007556675062972
$ws1.Cells.Item(65, 1).Value = '(''characters'', ''movie'')'
$ws1.Cells.Item(65, 2).Value = 4
$ws1.Cells.Item(65, 3).Value = 0.01007556675062972

$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 1).Value = '(''new'', ''york'')'
$ws2.Cells.Item(2, 2).Value = 12
$ws2.Cells.Item(2, 3).Value = 0.03934426229508197
$ws2.Cells.Item(3, 1).Value = '(''film'', ''also'')'
$ws2.Cells.Item(3, 2).Value = 9
$ws2.Cells.Item(3, 3).Value = 0.02950819672131148
$ws2.Cells.Item(4, 1).Value = '(''pretty'', ''good'')'
$ws2.Cells.Item(4, 2).Value = 9
$ws2.Cells.Item(4, 3).Value = 0.02950819672131148
$ws2.Cells.Item(5, 1).Value = '(''movie'', ''made'')'
$ws2.Cells.Item(5, 2).Value = 8
$ws2.Cells.Item(5, 3).Value = 0.02622950819672131
$ws2.Cells.Item(6, 1).Value = '(''rob'', ''roy'')'
$ws2.Cells.Item(6, 2).Value = 8
$ws2.Cells.Item(6, 3).Value = 0.02622950819672131
$ws2.Cells.Item(7, 1).Value = '(''quite'', ''good'')'
$ws2.Cells.Item(7, 2).Value = 7
$ws2.Cells.Item(7, 3).Value = 0.02295081967213115
$ws2.Cells.Item(8, 1).Value = '(''young'', ''man'')'
$ws2.Cells.Item(8, 2).Value = 6
$ws2.Cells.Item(8, 3).Value = 0.01967213114754099
$ws2.Cells.Item(9, 1).Value = '(''one'', ''people'')'
$ws2.Cells.Item(9, 2).Value = 6
$ws2.Cells.Item(9, 3).Value = 0.01967213114754099
$ws2.Cells.Item(10, 1).Value = '(''anna'', ''christie'')'
$ws2.Cells.Item(10, 2).Value = 6
$ws2.Cells.Item(10, 3).Value = 0.01967213114754099
$ws2.Cells.Item(11, 1).Value = '(''black'', ''white'')'
$ws2.Cells.Item(11, 2).Value = 6
$ws2.Cells.Item(11, 3).Value = 0.01967213114754099
$ws2.Cells.Item(12, 1).Value = '(''well'', ''done'')'
$ws2.Cells.Item(12, 2).Value = 5
$ws2.Cells.Item(12, 3).Value = 0.01639344262295082
$ws2.Cells.Item(13, 1).Value = '(''good'', ''film'')'
$ws2.Cells.Item(13, 2).Value = 5
$ws2.Cells.Item(13, 3).Value = 0.01639344262295082
$ws2.Cells.Item(14, 1).Value = '(''would'', ''like'')'
$ws2.Cells.Item(14, 2).Value = 5
$ws2.Cells.Item(14, 3).Value = 0.01639344262295082
$ws2.Cells.Item(15, 1).Value = '(''great'', ''deal'')'
$ws2.Cells.Item(15, 2).Value = 5
$ws2.Cells.Item(15, 3).Value = 0.01639344262295082
$ws2.Cells.Item(16, 1).Value = '(''great'', ''film'')'
$ws2.Cells.Item(16, 2).Value = 5
$ws2.Cells.Item(16, 3).Value = 0.01639344262295082
$ws2.Cells.Item(17, 1).Value = '(''uncle'', ''philip'')'
$ws2.Cells.Item(17, 2).Value = 5
$ws2.Cells.Item(17, 3).Value = 0.01639344262295082
$ws2.Cells.Item(18, 1).Value = '(''worth'', ''seeing'')'
$ws2.Cells.Item(18, 2).Value = 5
$ws2.Cells.Item(18, 3).Value = 0.01639344262295082
$ws2.Cells.Item(19, 1).Value = '(''film'', ''noir'')'
$ws2.Cells.Item(19, 2).Value = 5
$ws2.Cells.Item(19, 3).Value = 0.01639344262295082
$ws2.Cells.Item(20, 1).Value = '(''bug'', ''life'')'
$ws2.Cells.Item(20, 2).Value = 5
$ws2.Cells.Item(20, 3).Value = 0.01639344262295082
$ws2.Cells.Item(21, 1).Value = '(''melting'', ''man'')'
$ws2.Cells.Item(21, 2).Value = 5
$ws2.Cells.Item(21, 3).Value = 0.01639344262295082
$ws2.Cells.Item(22, 1).Value = '(''one'', ''day'')'
$ws2.Cells.Item(22, 2).Value = 5
$ws2.Cells.Item(22, 3).Value = 0.01639344262295082
$ws2.Cells.Item(23, 1).Value = '(''movie'', ''without'')'
$ws2.Cells.Item(23, 2).Value = 5
$ws2.Cells.Item(23, 3).Value = 0.01639344262295082
$ws2.Cells.Item(24, 1).Value = '(''real'', ''life'')'
$ws2.Cells.Item(24, 2).Value = 5
$ws2.Cells.Item(24, 3).Value = 0.01639344262295082
$ws2.Cells.Item(25, 1).Value = '(''well'', ''worth'')'
$ws2.Cells.Item(25, 2).Value = 4
$ws2.Cells.Item(25, 3).Value = 0.01311475409836066
$ws2.Cells.Item(26, 1).Value = '(''two'', ''films'')'
$ws2.Cells.Item(26, 2).Value = 4
$ws2.Cells.Item(26, 3).Value = 0.01311475409836066
$ws2.Cells.Item(27, 1).Value = '(''feel'', ''like'')'
$ws2.Cells.Item(27, 2).Value = 4
$ws2.Cells.Item(27, 3).Value = 0.01311475409836066
$ws2.Cells.Item(28, 1).Value = '(''get'', ''together'')'
$ws2.Cells.Item(28, 2).Value = 4
$ws2.Cells.Item(28, 3).Value = 0.01311475409836066
$ws2.Cells.Item(29, 1).Value = '(''takes'', ''place'')'
$ws2.Cells.Item(29, 2).Value = 4
$ws2.Cells.Item(29, 3).Value = 0.01311475409836066
$ws2.Cells.Item(30, 1).Value = '(''several'', ''times'')'
$ws2.Cells.Item(30, 2).Value = 4
$ws2.Cells.Item(30, 3).Value = 0.01311475409836066
$ws2.Cells.Item(31, 1).Value = '(''lion'', ''king'')'
$ws2.Cells.Item(31, 2).Value = 4
$ws2.Cells.Item(31, 3).Value = 0.01311475409836066
$ws2.Cells.Item(32, 1).Value = '(''well'', ''acted'')'
$ws2.Cells.Item(32, 2).Value = 4
$ws2.Cells.Item(32, 3).Value = 0.01311475409836066
$ws2.Cells.Item(33, 1).Value = '(''well'', ''written'')'
$ws2.Cells.Item(33, 2).Value = 4
$ws2.Cells.Item(33, 3).Value = 0.01311475409836066
$ws2.Cells.Item(34, 1).Value = '(''end'', ''well'')'
$ws2.Cells.Item(34, 2).Value = 4
$ws2.Cells.Item(34, 3).Value = 0.01311475409836066
$ws2.Cells.Item(35, 1).Value = '(''one'', ''liners'')'
$ws2.Cells.Item(35, 2).Value = 4
$ws2.Cells.Item(35, 3).Value = 0.01311475409836066
$ws2.Cells.Item(36, 1).Value = '(''well'', ''made'')'
$ws2.Cells.Item(36, 2).Value = 4
$ws2.Cells.Item(36, 3).Value = 0.01311475409836066
$ws2.Cells.Item(37, 1).Value = '(''car'', ''chases'')'
$ws2.Cells.Item(37, 2).Value = 4
$ws2.Cells.Item(37, 3).Value = 0.01311475409836066
$ws2.Cells.Item(38, 1).Value = '(''two'', ''years'')'
$ws2.Cells.Item(38, 2).Value = 4
$ws2.Cells.Item(38, 3).Value = 0.01311475409836066
$ws2.Cells.Item(39, 1).Value = '(''many'', ''ways'')'
$ws2.Cells.Item(39, 2).Value = 4
$ws2.Cells.Item(39, 3).Value = 0.01311475409836066
$ws2.Cells.Item(40, 1).Value = '(''b'', ''movie'')'
$ws2.Cells.Item(40, 2).Value = 4
$ws2.Cells.Item(40, 3).Value = 0.01311475409836066
$ws2.Cells.Item(41, 1).Value = '(''film'', ''festival'')'
$ws2.Cells.Item(41, 2).Value = 4
$ws2.Cells.Item(41, 3).Value = 0.01311475409836066
$ws2.Cells.Item(42, 1).Value = '(''muppet'', ''movie'')'
$ws2.Cells.Item(42, 2).Value = 4
$ws2.Cells.Item(42, 3).Value = 0.01311475409836066
$ws2.Cells.Item(43, 1).Value = '(''le'', ''million'')'
$ws2.Cells.Item(43, 2).Value = 4
$ws2.Cells.Item(43, 3).Value = 0.01311475409836066
$ws2.Cells.Item(44, 1).Value = '(''would'', ''never'')'
$ws2.Cells.Item(44, 2).Value = 4
$ws2.Cells.Item(44, 3).Value = 0.01311475409836066
$ws2.Cells.Item(45, 1).Value = '(''never'', ''get'')'
$ws2.Cells.Item(45, 2).Value = 4
$ws2.Cells.Item(45, 3).Value = 0.01311475409836066
$ws2.Cells.Item(46, 1).Value = '(''hit'', ''man'')'
$ws2.Cells.Item(46, 2).Value = 4
$ws2.Cells.Item(46, 3).Value = 0.01311475409836066
$ws2.Cells.Item(47, 1).Value = '(''never'', ''made'')'
$ws2.Cells.Item(47, 2).Value = 4
$ws2.Cells.Item(47, 3).Value = 0.01311475409836066
$ws2.Cells.Item(48, 1).Value = '(''art'', ''action'')'
$ws2.Cells.Item(48, 2).Value = 4
$ws2.Cells.Item(48, 3).Value = 0.01311475409836066
$ws2.Cells.Item(49, 1).Value = '(''beauty'', ''art'')'
$ws2.Cells.Item(49, 2).Value = 4
$ws2.Cells.Item(49, 3).Value = 0.01311475409836066
$ws2.Cells.Item(50, 1).Value = '(''stan'', ''ollie'')'
$ws2.Cells.Item(50, 2).Value = 4
$ws2.Cells.Item(50, 3).Value = 0.01311475409836066
$ws2.Cells.Item(51, 1).Value = '(''one'', ''would'')'
$ws2.Cells.Item(51, 2).Value = 4
$ws2.Cells.Item(51, 3).Value = 0.01311475409836066
$ws2.Cells.Item(52, 1).Value = '(''film'', ''really'')'
$ws2.Cells.Item(52, 2).Value = 4
$ws2.Cells.Item(52, 3).Value = 0.01311475409836066
$ws2.Cells.Item(53, 1).Value = '(''gives'', ''film'')'
$ws2.Cells.Item(53, 2).Value = 4
$ws2.Cells.Item(53, 3).Value = 0.01311475409836066
$ws2.Cells.Item(54, 1).Value = '(''people'', ''enjoy'')'
$ws2.Cells.Item(54, 2).Value = 4
$ws2.Cells.Item(54, 3).Value = 0.01311475409836066
$ws2.Cells.Item(55, 1).Value = '(''good'', ''little'')'
$ws2.Cells.Item(55, 2).Value = 4
$ws2.Cells.Item(55, 3).Value = 0.01311475409836066
$ws2.Cells.Item(56, 1).Value = '(''one'', ''favorite'')'
$ws2.Cells.Item(56, 2).Value = 4
$ws2.Cells.Item(56, 3).Value = 0.01311475409836066
$ws2.Cells.Item(57, 1).Value = '(''story'', ''great'')'
$ws2.Cells.Item(57, 2).Value = 4
$ws2.Cells.Item(57, 3).Value = 0.01311475409836066
$ws2.Cells.Item(58, 1).Value = '(''production'', ''values'')'
$ws2.Cells.Item(58, 2).Value = 4
$ws2.Cells.Item(58, 3).Value = 0.01311475409836066
$ws2.Cells.Item(59, 1).Value = '(''little'', ''film'')'
$ws2.Cells.Item(59, 2).Value = 4
$ws2.Cells.Item(59, 3).Value = 0.01311475409836066
$ws2.Cells.Item(60, 1).Value = '(''without'', ''doubt'')'
$ws2.Cells.Item(60, 2).Value = 4
$ws2.Cells.Item(60, 3).Value = 0.01311475409836066
$ws2.Cells.Item(61, 1).Value = '(''years'', ''later'')'
$ws2.Cells.Item(61, 2).Value = 4
$ws2.Cells.Item(61, 3).Value = 0.01311475409836066
$ws2.Cells.Item(62, 1).Value = '(''must'', ''see'')'
$ws2.Cells.Item(62, 2).Value = 4
$ws2.Cells.Item(62, 3).Value = 0.01311475409836066
$ws2.Cells.Item(63, 1).Value = '(''saw'', ''movie'')'
$ws2.Cells.Item(63, 2).Value = 4
$ws2.Cells.Item(63, 3).Value = 0.01311475409836066
$ws2.Cells.Item(64, 1).Value = '(''york'', ''love'')'
$ws2.Cells.Item(64, 2).Value = 4
$ws2.Cells.Item(64, 3).Value = 0.01311475409836066
$ws2.Cells.Item(65, 1).Value = '(''like'', ''tv'')'
$ws2.Cells.Item(65, 2).Value = 3
$ws2.Cells.Item(65, 3).Value = 0.009836065573770493

Write-Host "done"
